# Checklist.xlsx edit: "First Latex draft without code done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Task list (B11:D20): remove the old autofilter/data, rewrite the task
#    rows in their final order, then re-establish the autofilter on the
#    smaller B11:D14 range (it must be (re)created while rows 15-20 are
#    still empty, otherwise Excel auto-expands the filter to the whole
#    contiguous block).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("B12:N20").ClearContents()
$ws.Range("B11:D14").AutoFilter(1)

# Row 12
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = "Draft"
$ws.Range("D12").Value = "Add code to latex"

# Row 13
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Software - Logic"
$ws.Range("D13").Value = "Commenting"

# Row 14
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Software - Main"
$ws.Range("D14").Value = "Commenting"

# Row 15
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Software - Physical"
$ws.Range("D15").Value = "Commenting"

# Row 16
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Draft"
$ws.Range("D16").Value = "Add listings / Figures list"

# Row 19 written before row 17 so the two brand-new strings land in the
# shared-string table in the same order as the authored workbook
# ("Review Project Management section" before "Correct references from
# feedback").
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Draft"
$ws.Range("D19").Value = "Review Project Management section"

# Row 17
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Draft"
$ws.Range("D17").Value = "Correct references from feedback"

# Row 18
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Draft"
$ws.Range("D18").Value = "Format Pseudo-code figures"

# Row 20
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "Draft - Literature Review"
$ws.Range("D20").Value = "Add key technologies I will use"

# D15/D16 need to carry the same "no roadblock" fill as the rest of the
# column (previously they used the "waiting on implementation" fill) -
# copy the format from D12 so the existing style record is reused instead
# of a new one being synthesised.
$ws.Range("D12").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) New weekend plan rows (26:27)
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = "Saturday:"
$ws.Range("D27").Value = "Sunday:"
$ws.Range("E27").Value = "Final Edit"
$ws.Range("E26").Value = "Format Code"

# ---------------------------------------------------------------------------
# 3) Defined name _FilterDatabase must track the new autofilter range
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$B`$11:`$D`$14"
    }
}

# ---------------------------------------------------------------------------
# 4) Selection moves to D30 (last user action before saving)
# ---------------------------------------------------------------------------
$ws.Range("D30").Select()
